$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.990.00'
$ws.Range("E2").Value = '  +1.67%  '
$ws.Range("D3").Value = '3.097.20'
$ws.Range("E3").Value = '  +0.89%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''545.06'
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("D6").Value = '''140.68'
$ws.Range("E6").Value = '  +3.95%  '
$ws.Range("D7").Value = '''0.998'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '3.091.07'
$ws.Range("E8").Value = '  +0.87%  '
$ws.Range("D9").Value = '''0.499'
$ws.Range("E9").Value = '  +1.39%  '
$ws.Range("D10").Value = '''6.55'
$ws.Range("E10").Value = '  +3.68%  '
$ws.Range("D11").Value = '''0.157'
$ws.Range("E11").Value = '  +0.73%  '
$ws.Range("D12").Value = '''0.459'
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("E13").Value = '  +4.86%  '
$ws.Range("D14").Value = '''34.97'
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("D15").Value = '3.591.18'
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").Value = '64.034.94'
$ws.Range("E16").Value = '  +1.68%  '
$ws.Range("E17").Value = '  +1.42%  '
$ws.Range("D18").Value = '3.094.05'
$ws.Range("E18").Value = '  +0.79%  '
$ws.Range("D19").Value = '''6.68'
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("D20").Value = '''481.08'
$ws.Range("E20").Value = '  -1.09%  '
$ws.Range("D21").Value = '''13.48'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = '''0.702'
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("D23").Value = '''7.13'
$ws.Range("E23").Value = '  -0.90%  '
$ws.Range("D24").Value = '''79.10'
$ws.Range("E24").Value = '  +1.55%  '
$ws.Range("D25").Value = '''12.39'
$ws.Range("E25").Value = '  +0.72%  '
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("D27").Value = '''2.73'
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").Value = '''8.12'
$ws.Range("E28").Value = '  -4.08%  '
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("D30").Value = '''26.41'
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("E31").Value = '  -1.94%  '
$ws.Range("D32").Value = '''1.16'
$ws.Range("E32").Value = '  +2.77%  '
$ws.Range("D33").Value = '''57.37'
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("D34").Value = '''2.37'
$ws.Range("E34").Value = '  -5.22%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").Value = '''5.43'
$ws.Range("E35").Value = '  +6.24%  '
$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").Value = '''498.46'
$ws.Range("E36").Value = '  -3.72%  '
$ws.Range("D37").Value = '''6.04'
$ws.Range("E37").Value = '  +1.16%  '
$ws.Range("D38").Value = '3.256.98'
$ws.Range("E38").Value = '  +5.14%  '
$ws.Range("D39").Value = '''0.0406'
$ws.Range("E39").Value = '  +0.99%  '
$ws.Range("D40").Value = '''0.0804'
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("E41").Value = '  +1.06%  '
$ws.Range("D42").Value = '''2.72'
$ws.Range("E42").Value = '  +2.80%  '
$ws.Range("D43").Value = '''8.14'
$ws.Range("E43").Value = '  +0.72%  '
$ws.Range("D44").Value = '''0.255'
$ws.Range("E44").Value = '  +0.64%  '
$ws.Range("E46").Value = '  +2.94%  '
$ws.Range("D47").Value = '''25.39'
$ws.Range("E47").Value = '  +3.61%  '
$ws.Range("D48").Value = '''2.05'
$ws.Range("E48").Value = '  -1.23%  '
$ws.Range("D49").Value = '0.0₃0536'
$ws.Range("E49").Value = '  +7.90%  '
$ws.Range("D51").Value = '''2.41'
$ws.Range("E51").Value = '  +3.77%  '
